# Apply the "Create Batch_Miscellaneous sheet" edit to the CASH REPORT workbook.
#
# 1) BATCHID 274781 -> 915726 in Batch_Header (col B, rows 2-5) and
#    Batch_Detail (col C, rows 2-97).
# 2) Add a new "Batch_Miscellaneous" sheet (subset of the Batch_Header
#    columns) with a header row + 4 data rows, defaults set on a handful
#    of columns.

$wb = $excel.ActiveWorkbook

# --- 1) Update BATCHID values -------------------------------------------

$wsHeader = $wb.Worksheets.Item("Batch_Header")
for ($r = 2; $r -le 5; $r++) {
    $wsHeader.Cells.Item($r, 2).Value = 915726
}

$wsDetail = $wb.Worksheets.Item("Batch_Detail")
for ($r = 2; $r -le 97; $r++) {
    $wsDetail.Cells.Item($r, 3).Value = 915726
}

# --- 2) Create the Batch_Miscellaneous sheet ------------------------------

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsMisc = $wb.Worksheets.Add($null, $lastSheet)
$wsMisc.Name = "Batch_Miscellaneous"

$headers = @("BATCHID","ENTRYNO","DETAILNO","MISCCODE","NAME","ADDRESS1","ADDRESS2","ADDRESS3","ADDRESS4","POSTCODE","PHONE1","PHONE2","FAXNUMBER","CONTACT","COMMENTS","SWKEEPTOT","ACCTROW","ACCTNAME","ACCTNO","CITY","STATE","COUNTRY","EMAILADDR","URLADDR","SWAPPROVED","EFTDESC","BANKNAME","ACCOUNT","BRANCH","ACCTYPE","BANKID","SWIFTCTY","PAYDETL","ADDINFO1","ADDINFO2","COVERTYPE","COVERINFO","BENCODE","EITYPE","BOPCATG","BOPREF","BOPDESC","BRN","IDN")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $wsMisc.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# Batch_Miscellaneous DataFrame: a subset of Batch_Header (BATCHID, one row
# per Batch_Header detail entry), plus columns merged in from a second
# frame (DETAILNO) and a handful of explicit column defaults.
$batchId = 915726
for ($r = 2; $r -le 5; $r++) {
    $entryNo = "'{0:D5}" -f ($r - 1)

    $wsMisc.Cells.Item($r, 1).Value = $batchId      # BATCHID
    $wsMisc.Cells.Item($r, 2).Value = $entryNo      # ENTRYNO
    $wsMisc.Cells.Item($r, 3).Value = "'0000000200" # DETAILNO

    # Defaults for some columns on the Batch_Miscellaneous DataFrame
    $wsMisc.Cells.Item($r, 16).Value = "'FALSE"     # SWKEEPTOT
    $wsMisc.Cells.Item($r, 17).Value = 1            # ACCTROW
    $wsMisc.Cells.Item($r, 25).Value = "'FALSE"     # SWAPPROVED
    $wsMisc.Cells.Item($r, 30).Value = 0            # ACCTYPE
    $wsMisc.Cells.Item($r, 36).Value = 0            # COVERTYPE
    $wsMisc.Cells.Item($r, 39).Value = 0            # EITYPE
}
